# Appends 5 new numbered list paragraphs to the end of the document body
# (after the last existing paragraph "آب", before the sectPr), reproducing
# the target OOXML exactly -- including the per-run <w:rFonts w:hint="..."/>
# splitting -- via a raw WordprocessingML (Flat OPC) injection so that every
# run property matches the target precisely.

$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)

$xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>محاسبه‌</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>پ</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>چ</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ده</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>تابع</w:t></w:r><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>نما</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>یی</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>آلفا</w:t></w:r><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>قنطورس</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>انتخابات</w:t></w:r><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ر</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>است</w:t></w:r><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>جمهور</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:rtl/><w:lang w:val=`"en-US`"/></w:rPr><w:t>ی</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"en-US`" w:bidi=`"fa-IR`"/></w:rPr><w:t>0x55</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$endRange.InsertXML($xml)
